$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moves from K9 to E12 ---
$ws.Range("E12").Select()

# --- Row 1 (header): drop "Nomor Anggota" column, shift the remaining
#     headers one column left (C->B, D->C, E->D, F->E(No Telepon), G->F(Role),
#     H->G(Password)); old H1 is now empty. ---
$ws.Range("B1").Value = "Pangkat"
$ws.Range("C1").Value = "Kualifikasi"
$ws.Range("D1").Value = "Nrp"
$ws.Range("E1").Value = "No Telepon"
$ws.Range("F1").Value = "Role"
$ws.Range("G1").Value = "Password"
$ws.Range("H1").ClearContents()

# --- Row 2 ---
$ws.Range("B2").ClearFormats()
$ws.Range("B2").Value = "Mayor"
$ws.Range("C2").Value = "Check Pilot"
$ws.Range("D2").Value = 7788999187
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").Value = "082244862271"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "admin"
$ws.Range("G2").Value = "skadron"
$ws.Range("H2").ClearContents()

# --- Row 3 ---
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "Mayor"
$ws.Range("C3").Value = "Check Pilot"
$ws.Range("D3").Value = 7788999188
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").Value = "082244862271"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").ClearFormats()
$ws.Range("F3").Value = "komandan"
$ws.Range("G3").Value = "skadron"
$ws.Range("H3").ClearContents()

# --- Row 4 ---
$ws.Range("B4").ClearFormats()
$ws.Range("B4").Value = "Mayor"
$ws.Range("C4").Value = "Check Pilot"
$ws.Range("D4").Value = 7788999189
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").Value = "082244862271"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").ClearFormats()
$ws.Range("F4").Value = "anggota"
$ws.Range("G4").Value = "skadron"
$ws.Range("H4").ClearContents()
